$d = $word.ActiveDocument

# --- Edit 1: Update the "June 2024" amount cell from "900" to "74934" ---
# The cell's text is split across two runs ("9" and "00"). We locate the
# range of the first occurrence of "900" (the June 2024 row, not the
# Total row further down) using Find with no replacement, then assign
# new text directly to that range so the runs get merged into one.
$amountRange = $d.Content
$found = $amountRange.Find.Execute("900", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $amountRange.Text = "74934"
}

# --- Edit 2: Simplify the declaration paragraph ---
# The original paragraph text is split across three runs, with
# proofErr "gramStart"/"gramEnd" markers bracketing the word "receipts".
# Replacing the whole sentence collapses it into a single run and drops
# the proofErr markers, while keeping the text itself unchanged.
$declarationText = "I hereby declare that all bills, receipts and information provided by me to claim this benefit is accurate and correct. I understand that any wrong/inappropriate/incorrect submission of bills/receipt/information will call for an appropriate disciplinary action that Company may deem fit."
$d.Content.Find.Execute($declarationText, $true, $false, $false, $false, $false, $true, 1, $false, $declarationText, 2)
